# issue #5: stock data output to json file
#
# On the "股票" (stock) sheet:
#   1. Clean up a handful of company names that had a stray inserted
#      space / period baked into them.
#   2. Insert a new "property_category" column right before the existing
#      "date" column, and stamp every data row with the literal value
#      "stock" so the JSON export can tell which property type each row
#      belongs to.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# --- 1. Fix mangled company names (extra space / stray period) ---------
[void]$ws.Cells.Replace("國泰金融控股股份有限公 司", "國泰金融控股股份有限公司")
[void]$ws.Cells.Replace("台新金融控股股份有限公 司", "台新金融控股股份有限公司")
[void]$ws.Cells.Replace("台_積體電路製造股份有 限公司", "台_積體電路製造股份有限公司")
[void]$ws.Cells.Replace("兆豐票f金融股份有限妗 司 .", "兆豐票f金融股份有限妗司.")
[void]$ws.Cells.Replace("台灣積體電路製造股份有 限公司", "台灣積體電路製造股份有限公司")

# --- 2. Insert the new property_category column (H) --------------------
# This pushes the former H (date), I (legislator_name) and J
# (legislator_id) columns one slot to the right, becoming I, J, K.
$ws.Columns("H:H").Insert()

$ws.Range("H1").Value = "property_category"

$lastRow = 12
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = "stock"
}
